# Mise à jour de certains champs de Modules et de Professeurs
#
# Renames two of the header cells on the only worksheet:
#   C1: "Enseignant"       -> "Chef  Module"
#   D1: "Nombre d'heures"  -> "Composants"
# Widens the two columns that now hold the longer labels, and moves the
# saved cell selection from D2 to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell text -------------------------------------------------
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# --- Column widths ------------------------------------------------------
# Column C -> 35 characters, Column D -> 24.5703125 characters (stored
# OOXML "width" units). ColumnWidth takes the value in the same character
# units; the host applies a +5/6 padding/pixel-rounding step internally,
# so feeding it (target - 5/6) reproduces the requested stored width.
$ws.Columns.Item(3).ColumnWidth = 35 - 5/6
$ws.Columns.Item(4).ColumnWidth = 24.5703125 - 5/6

# --- Selection ------------------------------------------------------
[void]$ws.Range("E8").Select()
